# Course Choose: Add for student course choose function
#
# The "课程" (Course) sheet gains two pieces of per-course info that used to
# be crammed into a single free-text "备注(Remarks)" column:
#   - H column is repurposed into "人数(Studnums)" (number of students, numeric)
#   - a new I column "课时(Coursenums)" holds the course-hours count (numeric)
# A few class times in column D also grow from a single weekly slot to two
# weekly slots (e.g. "星期三第3,4节" -> "星期三第3,4节 星期五第1,2节").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "课程" (Course) sheet is the first sheet

# --- Column layout: widen D, and give the new I column the same width as H ---
$ws.Columns.Item(4).ColumnWidth = (201 / 7)
$ws.Columns.Item(9).ColumnWidth = (123 / 7)

# --- Header row ---
$ws.Cells.Item(1, 8).Value = "人数(Studnums)"
$ws.Cells.Item(1, 9).Value = "课时(Coursenums)"

# --- Row 2: 高等数学1 / 阶101 ---
$ws.Cells.Item(2, 4).Value = "星期三第3,4节 星期五第1,2节"
$ws.Cells.Item(2, 8).Value = 40
$ws.Cells.Item(2, 9).Value = 80

# --- Row 3: 高等数学1 / 阶102 ---
$ws.Cells.Item(3, 4).Value = "星期三第3,4节 星期五第1,2节"
$ws.Cells.Item(3, 8).Value = 40
$ws.Cells.Item(3, 9).Value = 80

# --- Row 4: 物理1 / 普110 ---
$ws.Cells.Item(4, 4).Value = "星期二第3,4节 星期四第1,2节"
$ws.Cells.Item(4, 8).Value = 60
$ws.Cells.Item(4, 9).Value = 80

# --- Row 5: 物理2 / 普111 ---
$ws.Cells.Item(5, 4).Value = "星期二第5,6节 星期四第7,8节"
$ws.Cells.Item(5, 8).Value = 60
$ws.Cells.Item(5, 9).Value = 100

# --- Row 12: C语言程序设计 -- old "60课时" Remarks becomes 60 Coursenums, no Studnums ---
$ws.Cells.Item(12, 8).ClearContents()
$ws.Cells.Item(12, 9).Value = 60

# --- Row 13: 大学生就业指导 -- old "4课时" Remarks becomes 4 Coursenums ---
$ws.Cells.Item(13, 8).ClearContents()
$ws.Cells.Item(13, 9).Value = 4

# --- Row 14: 生产实习 -- old "8课时" Remarks becomes 8 Coursenums ---
$ws.Cells.Item(14, 8).ClearContents()
$ws.Cells.Item(14, 9).Value = 8

# --- View: scroll so column C is leftmost, with G14 the active selection ---
$ws.Activate()
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 3
$ws.Range("G14").Select()
